# Adds two new validator-rule rows ("latitude_new" / "longitude_new") to the
# "fakedata" sheet, inserted above the existing "latitude" / "longitude" rows
# (old rows 9 & 10), shifting everything below down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above row 9 so the existing data (old rows 9-36) moves
# down to rows 11-38.
$ws.Rows.Item(9).Resize(2).Insert()

# Fill column-by-column (A9, A10, B9, B10, ..., F9, F10) so new shared-string
# entries land in the same order the source workbook uses.
$ws.Cells.Item(9, 1).Value = "latitude_new"
$ws.Cells.Item(10, 1).Value = "longitude_new"

$ws.Cells.Item(9, 2).Value = "The new latitude of the water system"
$ws.Cells.Item(10, 2).Value = "The new longitude of the water system"

$ws.Cells.Item(9, 4).Value = 34.460500000000003
$ws.Cells.Item(10, 4).Value = -117.3646

$ws.Cells.Item(9, 5).Value = "warning"
$ws.Cells.Item(10, 5).Value = "warning"

$ws.Cells.Item(9, 6).Value = 'grepl("^-?\\d+(\\.\\d{1,8})*$", latitude_new) & !is.na(latitude_new)'
$ws.Cells.Item(10, 6).Value = 'grepl("^-?\\d+(\\.\\d{1,8})*$", longitude_new) & !is.na(longitude_new)'

# Match the resulting selection shown in the saved file.
$ws.Range("F11").Select()
